$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-15 03:29:25"
$wsZh.Range("H4").Value = "2016-03-15 03:30:19"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-15 03:29:34"
$wsDe.Range("H4").Value = "2016-03-15 03:30:33"
